$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.771.66"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "2.216.24"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'241.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.02%  "
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").Value = "'72.26"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.93%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.593"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.01%  "
$ws.Range("D10").Value = "'41.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.86%  "
$ws.Range("D11").Value = "'0.0945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "'6.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.20%  "
$ws.Range("D14").Value = "2.549.10"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").Value = "'14.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.74%  "
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").Value = "2.208.98"
$ws.Range("E17").Value = "  -1.56%  "
$ws.Range("D18").Value = "41.657.94"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("D20").Value = "'6.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "'72.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'11.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +23.73%  "
$ws.Range("D23").Value = "'229.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.03%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'11.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D28").Value = "'2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").Value = "'167.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'20.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "'0.0795"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.42%  "
$ws.Range("D33").Value = "'5.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("D34").Value = "'29.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  -11.56%  "
$ws.Range("E37").Value = "  -7.57%  "
$ws.Range("E38").Value = "  -4.43%  "
$ws.Range("D39").Value = "'13.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.81%  "
$ws.Range("D40").Value = "'2.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.00%  "
$ws.Range("D41").Value = "'5.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "'102.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.53%  "
$ws.Range("D46").Value = "'0.0998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("D47").Value = "'2.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.78%  "
$ws.Range("D48").Value = "'1.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "2.422.57"
$ws.Range("E51").Value = "  -1.64%  "
